$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-09-19"

# Update the column header text (shared string) for column B
$ws.Range("B1").Value2 = "September 2022 (through September 19)"

# Update/insert carjacking counts for September data (year-over-year columns)
# B=current year(2022), K=2021, T=2020, AC=2019, AL=2018, AU=2017, BD=2016, BM=2015
$ws.Range("AC2").Value2 = 2
$ws.Range("AL2").Value2 = 6
$ws.Range("B2").Value2 = 5
$ws.Range("K2").Value2 = 8
$ws.Range("T2").Value2 = 8
$ws.Range("AC3").Value2 = 3
$ws.Range("B3").Value2 = 5
$ws.Range("BM3").Value2 = 2
$ws.Range("K3").Value2 = 13
$ws.Range("T3").Value2 = 5
$ws.Range("K5").Value2 = 2
$ws.Range("AL6").Value2 = 2
$ws.Range("K6").Value2 = 2
$ws.Range("T7").Value2 = 2
$ws.Range("K9").Value2 = 3
$ws.Range("AL10").Value2 = 2
$ws.Range("B10").Value2 = 4
$ws.Range("BM10").Value2 = 2
$ws.Range("K10").Value2 = 6
$ws.Range("AU12").Value2 = 1
$ws.Range("T12").Value2 = 6
$ws.Range("BD13").Value2 = 1
$ws.Range("K14").Value2 = 5
$ws.Range("BM15").Value2 = 1
$ws.Range("B18").Value2 = 1
$ws.Range("BM18").Value2 = 1
$ws.Range("B20").Value2 = 5
$ws.Range("AL21").Value2 = 2
$ws.Range("BD21").Value2 = 1
$ws.Range("B23").Value2 = 3
$ws.Range("K24").Value2 = 2
$ws.Range("T24").Value2 = 3
$ws.Range("B27").Value2 = 1
$ws.Range("BD27").Value2 = 1
$ws.Range("K27").Value2 = 1
$ws.Range("AC28").Value2 = 1
$ws.Range("AC33").Value2 = 1
$ws.Range("K33").Value2 = 2
$ws.Range("AC38").Value2 = 2
$ws.Range("K38").Value2 = 4
$ws.Range("K40").Value2 = 1
$ws.Range("K42").Value2 = 3
$ws.Range("AU44").Value2 = 1
$ws.Range("B44").Value2 = 1
$ws.Range("K47").Value2 = 1
$ws.Range("K51").Value2 = 1
$ws.Range("K74").Value2 = 1
$ws.Range("K77").Value2 = 3
$ws.Range("T94").Value2 = 1
$ws.Range("BD98").Value2 = 1
